# Apply "Transfer credit between accounts" updates to the expenditure register report.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column B ("Περίοδος") moves from period 9 to period 10 for every data row (2-20).
for ($row = 2; $row -le 20; $row++) {
    $ws.Cells.Item($row, 2).Value = 10
}

# Row 2 (ΑΛΕ 2120201017): invoice amount increases, new outstanding obligations of 230.
$ws.Range("M2").Value = 51460
$ws.Range("R2").Value = 230
$ws.Range("T2").Value = 230
$ws.Range("V2").Value = 230

# Row 4 (ΑΛΕ 2120211001): invoice amount increases, new outstanding obligations of 199.16.
$ws.Range("M4").Value = 4538.99
$ws.Range("R4").Value = 199.16
$ws.Range("T4").Value = 199.16
$ws.Range("V4").Value = 199.16

# Row 5 (ΑΛΕ 2120207001): invoice amount increases, new outstanding obligations of 800.
$ws.Range("M5").Value = 222070.15
$ws.Range("R5").Value = 800

# Row 8 (ΑΛΕ 2420403001): invoice amount increases, new outstanding obligations of 40.
$ws.Range("M8").Value = 3170
$ws.Range("R8").Value = 40
$ws.Range("T8").Value = 40
$ws.Range("V8").Value = 40

# Row 9 (ΑΛΕ 2420404001): invoice amount increases, new outstanding obligations of 1261.65.
$ws.Range("M9").Value = 30501.95
$ws.Range("R9").Value = 1261.65
$ws.Range("T9").Value = 1261.65
$ws.Range("V9").Value = 1261.65

# Row 13 (ΑΛΕ 2390501001): credit transferred between accounts -> invoice/payment up, commitments down.
$ws.Range("M13").Value = 7838.51
$ws.Range("N13").Value = 7838.51
$ws.Range("Q13").Value = 236.49

# Row 18 (ΑΛΕ 2390589001): credit transferred between accounts -> invoice/payment up, commitments down to 0.
$ws.Range("M18").Value = 4576
$ws.Range("N18").Value = 4576
$ws.Range("Q18").Value = 0

# Row 20 (ΑΛΕ 2420407001): new invoice amount and matching outstanding obligations of 507.99.
$ws.Range("M20").Value = 507.99
$ws.Range("R20").Value = 507.99
$ws.Range("T20").Value = 507.99
$ws.Range("V20").Value = 507.99
